# Clear the "额外的奖励" (ExtendReward) values in column E for rows 5-11.
# Only the values are removed; number formatting / style stays untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5:E11").ClearContents()

# Update the active selection to reflect the latest manual edit location (C5).
$ws.Range("C5").Select()
